$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(38).Insert()

$ws.Cells.Item(38, 1).Value = 3
$ws.Cells.Item(38, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(38, 3).Value = "Coquimbo"
$ws.Cells.Item(38, 4).Value = 44494
$ws.Cells.Item(38, 5).Value = 5
$ws.Cells.Item(38, 6).Value = 100112026
$ws.Cells.Item(38, 7).Value = "Haba"
$ws.Cells.Item(38, 8).Value = "Sin especificar"
$ws.Cells.Item(38, 9).Value = "Primera"
$ws.Cells.Item(38, 10).Value = 73
$ws.Cells.Item(38, 11).Value = 8000
$ws.Cells.Item(38, 12).Value = 8500
$ws.Cells.Item(38, 13).Value = 8240
$ws.Cells.Item(38, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(38, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(38, 16).Value = 330
$ws.Cells.Item(38, 17).Value = 25
$ws.Cells.Item(38, 18).Value = "Hortaliza"
